# Apply updates to "Teams Data" sheet:
# - Update Powerups (column H) for the IMPOSTORS team row (row 24)
# - Update Users (column C) for the Chdi Gang team row (row 25)
# - Update Score (column J) values for rows 2, 3, 7, 24

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teams Data")

# Update power-ups data structure for team "IMPOSTORS" (row 24)
$ws.Range("H24").Value = "2, 4, 7, 6"

# Update task data / team member description for team "Chdi Gang" (row 25)
$ws.Range("C25").Value = "Sameer Verma, Aditya, Palak, Bhavninder"

# Update Score values
$ws.Range("J2").Value = 2985
$ws.Range("J3").Value = 1618
$ws.Range("J7").Value = 7282
$ws.Range("J24").Value = 454

$wb.Save()
